$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to Text format first, otherwise Excel will auto-convert them to
# numeric values (losing the intended display format, e.g. trailing zeros).
$textForcedRefs = @('D5', 'D6', 'D7', 'D10', 'D11', 'D12', 'D16', 'D18', 'D21', 'D23', 'D24', 'D25', 'D28', 'D30', 'D31', 'D32', 'D35', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D45', 'D46', 'D47', 'D51')
foreach ($ref in $textForcedRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated coin names / links / prices / volume percentages.
$ws.Range('D2').Value = '42.457.09'
$ws.Range('E2').Value = '  -2.47%  '
$ws.Range('D3').Value = '2.275.17'
$ws.Range('E3').Value = '  -4.36%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '300.04'
$ws.Range('E5').Value = '  -3.24%  '
$ws.Range('D6').Value = '96.81'
$ws.Range('E6').Value = '  -7.34%  '
$ws.Range('D7').Value = '0.502'
$ws.Range('E7').Value = '  -1.98%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -5.24%  '
$ws.Range('D10').Value = '33.57'
$ws.Range('E10').Value = '  -6.68%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = '0.0786'
$ws.Range('E11').Value = '  -3.23%  '
$ws.Range('B12').Value = 'OKB'
$ws.Range('C12').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D12').Value = '50.55'
$ws.Range('E12').Value = '  -5.03%  '
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('E14').Value = '  -4.91%  '
$ws.Range('D15').Value = '2.628.68'
$ws.Range('E15').Value = '  -4.42%  '
$ws.Range('D16').Value = '15.12'
$ws.Range('E16').Value = '  -3.23%  '
$ws.Range('D17').Value = '2.261.68'
$ws.Range('E17').Value = '  -5.09%  '
$ws.Range('D18').Value = '0.783'
$ws.Range('E18').Value = '  -3.47%  '
$ws.Range('D19').Value = '42.370.28'
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('E20').Value = '  -2.81%  '
$ws.Range('D21').Value = '11.37'
$ws.Range('E21').Value = '  -4.46%  '
$ws.Range('E22').Value = '  -5.29%  '
$ws.Range('D23').Value = '66.47'
$ws.Range('E23').Value = '  -2.77%  '
$ws.Range('D24').Value = '234.97'
$ws.Range('E24').Value = '  -2.32%  '
$ws.Range('D25').Value = '1.92'
$ws.Range('E25').Value = '  -6.10%  '
$ws.Range('E26').Value = '  -5.28%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = '24.32'
$ws.Range('E28').Value = '  -5.81%  '
$ws.Range('E29').Value = '  -1.83%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '164.51'
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '33.65'
$ws.Range('E31').Value = '  -8.25%  '
$ws.Range('D32').Value = '9.05'
$ws.Range('E32').Value = '  -4.84%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('E34').Value = '  -6.19%  '
$ws.Range('D35').Value = '2.40'
$ws.Range('E35').Value = '  -4.91%  '
$ws.Range('E36').Value = '  -6.05%  '
$ws.Range('D37').Value = '4.31'
$ws.Range('E37').Value = '  -8.52%  '
$ws.Range('D38').Value = '2.81'
$ws.Range('E38').Value = '  -9.31%  '
$ws.Range('D39').Value = '16.06'
$ws.Range('E39').Value = '  -12.08%  '
$ws.Range('D40').Value = '0.0996'
$ws.Range('E40').Value = '  -6.04%  '
$ws.Range('D41').Value = '1.75'
$ws.Range('E41').Value = '  -9.20%  '
$ws.Range('D42').Value = '0.109'
$ws.Range('E42').Value = '  -3.77%  '
$ws.Range('E43').Value = '  -7.54%  '
$ws.Range('D44').Value = '1.955.58'
$ws.Range('E44').Value = '  -3.90%  '
$ws.Range('D45').Value = '0.0281'
$ws.Range('E45').Value = '  -3.63%  '
$ws.Range('D46').Value = '17.85'
$ws.Range('E46').Value = '  -9.67%  '
$ws.Range('D47').Value = '9.65'
$ws.Range('E47').Value = '  -8.68%  '
$ws.Range('E48').Value = '  -10.10%  '
$ws.Range('E49').Value = '  -4.64%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.500.29'
$ws.Range('E50').Value = '  -4.17%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '4.66'
$ws.Range('E51').Value = '  -1.55%  '
